$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting existing rows 126:200 down to 127:201.
$ws.Rows(126).Insert()

# Populate the newly inserted row 126 with its data (mirrors the surrounding rows
# for the fixed market/product columns, plus the new record's own values).
$ws.Range("A126").Value = 4
$ws.Range("B126").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C126").Value = "Los Lagos"
$ws.Range("D126").Value = 44777
$ws.Range("E126").Value = 10
$ws.Range("F126").Value = 100112009
$ws.Range("G126").Value = "Acelga"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 200
$ws.Range("K126").Value = 1200
$ws.Range("L126").Value = 1500
$ws.Range("M126").Value = 1350
$ws.Range("N126").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O126").Value = "Región de Los Lagos"
$ws.Range("P126").Value = 900
$ws.Range("Q126").Value = 1.5
$ws.Range("R126").Value = "Hortaliza"
